$d = $word.ActiveDocument

# --- Hunk 1: merge " (a " + "clique" + ")." into a single run ---
$d.Content.Find.Execute(
    " (a clique).", $true, $false, $false, $false, $false,
    $true, 1, $false, " (a clique).", 2) | Out-Null

# --- Hunk 2: merge "dynamic" ... " between runs." into a single run ---
$dash = [char]0x2013
$hunk2Text = "dynamic " + $dash + " every node is connected to 3 nodes randomly and another connection to the node with the next index (to ensure full connectivity). Not consistent between runs."
$d.Content.Find.Execute(
    $hunk2Text,
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    $hunk2Text,
    2) | Out-Null

# --- Hunk 3: merge "python3.7 " + "make_setup_test.py" + " 2 200 20 mesh" ---
$d.Content.Find.Execute(
    "python3.7 make_setup_test.py 2 200 20 mesh", $true, $false, $false, $false, $false,
    $true, 1, $false, "python3.7 make_setup_test.py 2 200 20 mesh", 2) | Out-Null

# --- Hunk 5: merge "Range of private IP addresses, if different from " + "AWS setup requirements" + "." ---
$d.Content.Find.Execute(
    "Range of private IP addresses, if different from AWS setup requirements.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Range of private IP addresses, if different from AWS setup requirements.", 2) | Out-Null

# --- Hunk 4: relocate the _GoBack bookmark and highlight "run_multiple_tests" ---

# Highlight just the "run_multiple_tests" substring (leaving ".py, that can run multiple" as plain text)
$hl = $d.Content
$hl.Find.Execute("run_multiple_tests", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$hl.Font.HighlightColorIndex = 7

# Move the "_GoBack" bookmark from the end of the "...topology of a clique." paragraph
# to the middle of "according" ("acco" | "rding to lists"), matching where the edit
# actually left off. Adding a bookmark with an existing name replaces the old one.
$bm = $d.Content
$bm.Find.Execute("according to lists", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bm.Start = $bm.Start + 4
$bm.End = $bm.Start
$d.Bookmarks.Add("_GoBack", $bm) | Out-Null
